# Auto-generated edit script: applies the row-content rotation
# described by the commit diff to rows 114-127 of sheet "Artfynd".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 114
$ws.Range("A114").Value2 = 111743546
$ws.Range("Q114").Value2 = 339474.5644867857
$ws.Range("R114").Value2 = 6571113.931964876

# Row 115
$ws.Range("A115").Value2 = 111743520
$ws.Range("B115").Value2 = 56398
$ws.Range("D115").Value2 = "NT"
$ws.Range("E115").Value2 = 100109
$ws.Range("F115").Value2 = "Tretåig hackspett"
$ws.Range("G115").Value2 = "Picoides tridactylus"
$ws.Range("H115").Value2 = "(Linnaeus, 1758)"
$ws.Range("M115").Value2 = "färska spår"
$ws.Range("Q115").Value2 = 339096.8530521042
$ws.Range("R115").Value2 = 6571013.66294401
$ws.Range("AJ115").Value2 = "gran"
$ws.Range("AK115").Value2 = "Picea abies"
$ws.Range("AO115").Value2 = "Picea abies"

# Row 117
$ws.Range("A117").Value2 = 111743519
$ws.Range("B117").Value2 = 90666
$ws.Range("D117").Value2 = "LC"
$ws.Range("E117").Value2 = 4364
$ws.Range("F117").Value2 = "Dropptaggsvamp"
$ws.Range("G117").Value2 = "Hydnellum ferrugineum"
$ws.Range("H117").Value2 = "(Fr.:Fr.) P. Karst."
$ws.Range("Q117").Value2 = 339118.4126724883
$ws.Range("R117").Value2 = 6571062.424656671

# Row 118
$ws.Range("A118").Value2 = 111743526
$ws.Range("B118").Value2 = 90666
$ws.Range("D118").Value2 = "LC"
$ws.Range("E118").Value2 = 4364
$ws.Range("F118").Value2 = "Dropptaggsvamp"
$ws.Range("G118").Value2 = "Hydnellum ferrugineum"
$ws.Range("H118").Value2 = "(Fr.:Fr.) P. Karst."
$ws.Range("Q118").Value2 = 338870.1217119552
$ws.Range("R118").Value2 = 6571086.774471543

# Row 119
$ws.Range("A119").Value2 = 111743551
$ws.Range("B119").Value2 = 96348
$ws.Range("D119").Value2 = "VU"
$ws.Range("E119").Value2 = 220787
$ws.Range("F119").Value2 = "Knärot"
$ws.Range("G119").Value2 = "Goodyera repens"
$ws.Range("H119").Value2 = "(L.) R. Br."
$ws.Range("Q119").Value2 = 339522.8608171764
$ws.Range("R119").Value2 = 6571091.407599592
$ws.Range("AJ119").ClearContents()
$ws.Range("AK119").ClearContents()
$ws.Range("AO119").ClearContents()

# Row 120
$ws.Range("A120").Value2 = 111743554
$ws.Range("B120").Value2 = 88966
$ws.Range("D120").Value2 = "NT"
$ws.Range("E120").Value2 = 5754
$ws.Range("F120").Value2 = "Gultoppig fingersvamp"
$ws.Range("G120").Value2 = "Ramaria testaceoflava"
$ws.Range("H120").Value2 = "(Bres.) Corner"
$ws.Range("Q120").Value2 = 339577.2032005055
$ws.Range("R120").Value2 = 6571127.007499221

# Row 121
$ws.Range("A121").Value2 = 111743523
$ws.Range("B121").Value2 = 73634
$ws.Range("D121").Value2 = "LC"
$ws.Range("E121").Value2 = 6426
$ws.Range("F121").Value2 = "Kattfotslav"
$ws.Range("G121").Value2 = "Felipes leucopellaeus"
$ws.Range("H121").Value2 = "(Ach.) Frisch & G.Thor"
$ws.Range("J121").ClearContents()
$ws.Range("K121").ClearContents()
$ws.Range("L121").ClearContents()
$ws.Range("N121").ClearContents()
$ws.Range("Q121").Value2 = 339009.0243061834
$ws.Range("R121").Value2 = 6571011.238422027
$ws.Range("AF121").ClearContents()

# Row 122
$ws.Range("A122").Value2 = 111743549
$ws.Range("J122").Value2 = ""
$ws.Range("K122").Value2 = "blomning"
$ws.Range("L122").Value2 = ""
$ws.Range("N122").Value2 = ""
$ws.Range("Q122").Value2 = 339495.029088294
$ws.Range("R122").Value2 = 6571076.196190646
$ws.Range("AF122").Value2 = ""

# Row 123
$ws.Range("A123").Value2 = 111743517
$ws.Range("B123").Value2 = 73634
$ws.Range("D123").Value2 = "LC"
$ws.Range("E123").Value2 = 6426
$ws.Range("F123").Value2 = "Kattfotslav"
$ws.Range("G123").Value2 = "Felipes leucopellaeus"
$ws.Range("H123").Value2 = "(Ach.) Frisch & G.Thor"
$ws.Range("Q123").Value2 = 339278.3213300391
$ws.Range("R123").Value2 = 6571107.378548244

# Row 124
$ws.Range("A124").Value2 = 111743521
$ws.Range("B124").Value2 = 96348
$ws.Range("D124").Value2 = "VU"
$ws.Range("E124").Value2 = 220787
$ws.Range("F124").Value2 = "Knärot"
$ws.Range("G124").Value2 = "Goodyera repens"
$ws.Range("H124").Value2 = "(L.) R. Br."
$ws.Range("Q124").Value2 = 339070.1946752003
$ws.Range("R124").Value2 = 6571001.989220584

# Row 125
$ws.Range("A125").Value2 = 111743515
$ws.Range("B125").Value2 = 96348
$ws.Range("D125").Value2 = "VU"
$ws.Range("E125").Value2 = 220787
$ws.Range("F125").Value2 = "Knärot"
$ws.Range("G125").Value2 = "Goodyera repens"
$ws.Range("H125").Value2 = "(L.) R. Br."
$ws.Range("M125").ClearContents()
$ws.Range("Q125").Value2 = 339441.7613444271
$ws.Range("R125").Value2 = 6571017.506567059
$ws.Range("AJ125").ClearContents()
$ws.Range("AK125").ClearContents()
$ws.Range("AO125").ClearContents()

# Row 126
$ws.Range("A126").Value2 = 111743524
$ws.Range("B126").Value2 = 94134
$ws.Range("D126").Value2 = "NT"
$ws.Range("E126").Value2 = 53
$ws.Range("F126").Value2 = "Vedtrappmossa"
$ws.Range("G126").Value2 = "Crossocalyx hellerianus"
$ws.Range("H126").Value2 = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q126").Value2 = 338949.7235384365
$ws.Range("R126").Value2 = 6571040.381812023
$ws.Range("AJ126").Value2 = "tall"
$ws.Range("AK126").Value2 = "Pinus sylvestris"
$ws.Range("AO126").Value2 = "Pinus sylvestris"

# Row 127
$ws.Range("A127").Value2 = 111743516
$ws.Range("B127").Value2 = 96348
$ws.Range("D127").Value2 = "VU"
$ws.Range("E127").Value2 = 220787
$ws.Range("F127").Value2 = "Knärot"
$ws.Range("G127").Value2 = "Goodyera repens"
$ws.Range("H127").Value2 = "(L.) R. Br."
$ws.Range("Q127").Value2 = 339415.5147437509
$ws.Range("R127").Value2 = 6571015.54325202

